# The workbook's single sheet ("Artfynd") holds one species-observation
# per row (rows 2-12, columns A:AY). The edit re-shuffles which row each
# observation's data lives in (two permutation cycles: 2<->7, and
# 3->9->8->6->5->4->3, and 10->12->11->10) while every other row
# (header row 1, sheet structure, etc.) stays put.
#
# Strategy: snapshot every source row's full A:AY values first (so the
# cyclic re-assignment can't clobber data it still needs to read later),
# then write each snapshot into its new row according to the mapping.
#
# A couple of text columns hold values that LOOK like dates/times/numbers
# ("2023-08-28", "00:00", "25", ...) but are stored as plain text in the
# file. Assigning such a string straight into a General-formatted cell
# makes Excel "smart convert" it into a real date/number, which would
# silently change the cell's type. To avoid that we temporarily mark
# those columns as Text ("@") before writing, then clear the formatting
# back off again afterwards so the file doesn't end up with stray
# number-format differences on those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 12

# Columns whose text values could be misread as dates/times/numbers.
$textColumns = @("I", "Y", "Z", "AA", "AB")
foreach ($col in $textColumns) {
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).NumberFormat = "@"
}

# Snapshot every row (A:AY) before any writes happen.
$rowData = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData[$r] = $ws.Range("A" + $r + ":AY" + $r).Value2
}

# new row -> source row it should take its contents from.
$mapping = @{
    2  = 7
    3  = 9
    4  = 3
    5  = 4
    6  = 5
    7  = 2
    8  = 6
    9  = 8
    10 = 12
    11 = 10
    12 = 11
}

foreach ($newRow in $mapping.Keys) {
    $sourceRow = $mapping[$newRow]
    $ws.Range("A" + $newRow + ":AY" + $newRow).Value2 = $rowData[$sourceRow]
}

# Undo the temporary Text formatting so the cells end up with no
# explicit number format, same as before the script ran.
foreach ($col in $textColumns) {
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).ClearFormats()
}
